$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (target stored width 22.42578125 chars; the host
# quantizes ColumnWidth to 1/6-character increments, so 21.666666 is the
# input that lands closest/exactly on the nearest reachable stored width)
$ws.Columns.Item(1).ColumnWidth = 21.666666

# Row 1 - Date
$ws.Range("B1").Value = 44515
$ws.Range("C1").Value = $null

# Row 2 - Current Hashrate
$ws.Range("B2").Value = 178.9
$ws.Range("C2").Value = $null

# Row 3 - Unpaid ETH
$ws.Range("B3").Value = 0.0322
$ws.Range("C3").Value = $null

# Row 4 - Daily ETH
$ws.Range("B4").Value = 0.00296
$ws.Range("C4").Value = $null

# Row 5 - Sum of payouts
$ws.Range("C5").Value = $null

# Row 6 - Days to next payout
$ws.Range("B6").Value = 22
$ws.Range("C6").Value = $null

# Row 7 - new: Today's ETH price [USD]
$ws.Range("A7").Value = "Today's ETH price [USD]:"
$ws.Range("B7").Value = 4609.67

# Selection
$ws.Range("H9").Select()
